$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace each "NA" placeholder in column D with its own unique value,
# in the order the cells appear on the sheet (D5, D13, D17, D18, D19, D22).
$ws.Range("D5").Value = "removed_1"
$ws.Range("D13").Value = "removed_2"
$ws.Range("D17").Value = "removed_3"
$ws.Range("D18").Value = "removed_4"
$ws.Range("D19").Value = "removed_5"
$ws.Range("D22").Value = "removed_6"

# Move the active selection to D23
$ws.Range("D23").Select()
